$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'246.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Value = "'5"
$ws.Range("G2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'22.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").Value = "'5"
$ws.Range("G3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.435"
$ws.Range("D4").Style = "Normal"
$ws.Range("G4").Value = "'5"
$ws.Range("G4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.05786"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Value = "'5"
$ws.Range("G5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'3.393"
$ws.Range("D6").Style = "Normal"
$ws.Range("G6").Value = "'5"
$ws.Range("G6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'6.325"
$ws.Range("D7").Style = "Normal"
$ws.Range("G7").Value = "'5"
$ws.Range("G7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.8180"
$ws.Range("D8").Style = "Normal"
$ws.Range("G8").Value = "'5"
$ws.Range("G8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.9634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("G9").Value = "'5"
$ws.Range("G9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.1429"
$ws.Range("D10").Style = "Normal"
$ws.Range("G10").Value = "'5"
$ws.Range("G10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07493"
$ws.Range("D11").Style = "Normal"
$ws.Range("G11").Value = "'5"
$ws.Range("G11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.03147"
$ws.Range("D12").Style = "Normal"
$ws.Range("G12").Value = "'5"
$ws.Range("G12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.02993"
$ws.Range("D13").Style = "Normal"
$ws.Range("G13").Value = "'5"
$ws.Range("G13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'4.144"
$ws.Range("D14").Style = "Normal"
$ws.Range("G14").Value = "'5"
$ws.Range("G14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.09412"
$ws.Range("D15").Style = "Normal"
$ws.Range("G15").Value = "'5"
$ws.Range("G15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.001592"
$ws.Range("D16").Style = "Normal"
$ws.Range("G16").Value = "'5"
$ws.Range("G16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.04812"
$ws.Range("D17").Style = "Normal"
$ws.Range("G17").Value = "'5"
$ws.Range("G17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'0.0005849"
$ws.Range("D18").Style = "Normal"
$ws.Range("G18").Value = "'5"
$ws.Range("G18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.006184"
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").Value = "'5"
$ws.Range("G19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'0.004127"
$ws.Range("D20").Style = "Normal"
$ws.Range("G20").Value = "'5"
$ws.Range("G20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'0.0009948"
$ws.Range("D21").Style = "Normal"
$ws.Range("G21").Value = "'5"
$ws.Range("G21").Style = "Normal"

# Row 22
$ws.Range("G22").Value = "'5"
$ws.Range("G22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'3.772"
$ws.Range("D23").Style = "Normal"
$ws.Range("G23").Value = "'5"
$ws.Range("G23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'2.217"
$ws.Range("D24").Style = "Normal"
$ws.Range("G24").Value = "'5"
$ws.Range("G24").Style = "Normal"

# Row 25
$ws.Range("G25").Value = "'5"
$ws.Range("G25").Style = "Normal"

# Row 26
$ws.Range("G26").Value = "'5"
$ws.Range("G26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.0003998"
$ws.Range("D27").Style = "Normal"
$ws.Range("G27").Value = "'5"
$ws.Range("G27").Style = "Normal"

# Row 28
$ws.Range("G28").Value = "'5"
$ws.Range("G28").Style = "Normal"

# Row 29
$ws.Range("G29").Value = "'5"
$ws.Range("G29").Style = "Normal"

# Row 30
$ws.Range("G30").Value = "'5"
$ws.Range("G30").Style = "Normal"

# Row 31
$ws.Range("G31").Value = "'5"
$ws.Range("G31").Style = "Normal"

# Row 32
$ws.Range("G32").Value = "'5"
$ws.Range("G32").Style = "Normal"

# Row 33
$ws.Range("G33").Value = "'5"
$ws.Range("G33").Style = "Normal"

# Row 34
$ws.Range("G34").Value = "'5"
$ws.Range("G34").Style = "Normal"

# Row 35
$ws.Range("G35").Value = "'5"
$ws.Range("G35").Style = "Normal"

# Row 36
$ws.Range("G36").Value = "'5"
$ws.Range("G36").Style = "Normal"

# Row 37
$ws.Range("G37").Value = "'5"
$ws.Range("G37").Style = "Normal"

# Row 38
$ws.Range("G38").Value = "'5"
$ws.Range("G38").Style = "Normal"

# Row 39
$ws.Range("G39").Value = "'5"
$ws.Range("G39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.03891"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Value = "'5"
$ws.Range("G40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006458"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "'5"
$ws.Range("G41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1077"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "'5"
$ws.Range("G42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("G43").Value = "'5"
$ws.Range("G43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.006698"
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").Value = "'5"
$ws.Range("G44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.00005591"
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").Value = "'5"
$ws.Range("G45").Style = "Normal"

# Row 46
$ws.Range("G46").Value = "'5"
$ws.Range("G46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.3799"
$ws.Range("D47").Style = "Normal"
$ws.Range("G47").Value = "'5"
$ws.Range("G47").Style = "Normal"

# Row 48
$ws.Range("G48").Value = "'5"
$ws.Range("G48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").Value = "'5"
$ws.Range("G49").Style = "Normal"

# Row 50
$ws.Range("G50").Value = "'5"
$ws.Range("G50").Style = "Normal"

# Row 51
$ws.Range("G51").Value = "'5"
$ws.Range("G51").Style = "Normal"
